$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from an existing header cell (H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Keep the text values we just set (PasteSpecial formats only touches formatting)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 and IF data values per row (2-32)
$i0 = @{
    2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1; 11=1;
    12=6; 13=7; 14=6; 15=6; 16=6; 17=6;
    18=1; 19=1; 20=1; 21=1; 22=1; 23=1; 24=1; 25=1; 26=1; 27=1; 28=1; 29=1; 30=1;
    31=4; 32=1
}
$if = @{
    2=7; 3=4; 4=6; 5=6; 6=5; 7=7; 8=5; 9=4; 10=4; 11=3;
    12=8; 13=8; 14=6; 15=7; 16=7; 17=6;
    18=6; 19=4; 20=6; 21=5; 22=4; 23=5; 24=5; 25=6; 26=6; 27=5; 28=5; 29=4; 30=4;
    31=5; 32=2
}

for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 9).Value = $i0[$r]
    $ws.Cells.Item($r, 10).Value = $if[$r]
}
